# Added approve replenish page
# - Restocks the medicine in row 3 (Ibuprofen) to 50 units
# - Appends two new medicines: "CC7" (row 5) and "CCure" (row 6)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replenish existing stock level
$ws.Range("B3").Value = 50

# New medicine row: CCure (entered first so it lands earlier in the
# shared-string table, matching how the sheet was actually authored)
$ws.Range("A6").Value = "CCure"
$ws.Range("B6").Value = 470
$ws.Range("C6").Value = 10

# New medicine row: CC7
$ws.Range("A5").Value = "CC7"
$ws.Range("B5").Value = 200
$ws.Range("C5").Value = 20

# Leave the selection where the user ended up after the edits
$ws.Range("B7").Select()
